$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.304.60"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "3.388.22"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'252.69"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").Value = "'660.04"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("D8").Value = "'0.425"
$ws.Range("E8").Value = "  -2.85%  "

$ws.Range("D9").Value = "'1.04"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "3.385.37"
$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("E12").Value = "  -3.13%  "

$ws.Range("D13").Value = "'42.55"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "97.832.37"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "'6.13"
$ws.Range("E15").Value = "  -6.03%  "

$ws.Range("E16").Value = "  -4.30%  "

$ws.Range("D17").Value = "4.020.52"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").Value = "'9.18"
$ws.Range("E18").Value = "  +2.78%  "

$ws.Range("D19").Value = "3.382.48"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").Value = "'18.05"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").Value = "'0.520"

$ws.Range("D22").Value = "'11.01"
$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("D23").Value = "'511.85"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'3.43"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000202"
$ws.Range("E25").Value = "  -2.67%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'6.94"
$ws.Range("E26").Value = "  +2.83%  "

$ws.Range("D27").Value = "'96.85"
$ws.Range("E27").Value = "  -4.28%  "

$ws.Range("D28").Value = "'12.38"
$ws.Range("E28").Value = "  -4.42%  "

$ws.Range("D29").Value = "3.568.09"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "'11.80"
$ws.Range("E30").Value = "  +1.21%  "

$ws.Range("E31").Value = "  -3.90%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("E33").Value = "  -4.27%  "

$ws.Range("D34").Value = "'2.61"
$ws.Range("E34").Value = "  +8.52%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").Value = "'0.562"
$ws.Range("E36").Value = "  -2.29%  "

$ws.Range("D37").Value = "'28.92"
$ws.Range("E37").Value = "  -3.61%  "

$ws.Range("D38").Value = "'7.96"
$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("D39").Value = "'1.47"
$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("D40").Value = "'531.18"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("D41").Value = "'0.152"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("D44").Value = "'0.850"
$ws.Range("E44").Value = "  -3.22%  "

$ws.Range("D45").Value = "'0.0428"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").Value = "'1.73"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").Value = "'3.69"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.26"
$ws.Range("E48").Value = "  +5.65%  "

$ws.Range("D49").Value = "'8.66"
$ws.Range("E49").Value = "  -3.79%  "

$ws.Range("D50").Value = "'55.95"
$ws.Range("E50").Value = "  +3.75%  "

$ws.Range("D51").Value = "'5.58"
$ws.Range("E51").Value = "  -4.50%  "
